# Updated cryptos list (price + 1h volume change columns) with a GitHub Actions
# scrape refresh; also corrects the ranking swap between ImmutableX and
# PolygonEcosystemToken (rows 34/35).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must remain TEXT (e.g. "1.00", "545.15").
# Temporarily force Text format while assigning so Excel does not coerce them to numbers,
# then restore the original General format.
$dFormat = $ws.Range("D2").NumberFormat
$ws.Columns("D").NumberFormat = "@"

$ws.Range("D2").Value = "60.185.31"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").Value = "2.318.85"
$ws.Range("E3").Value = "  -1.42%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "545.15"
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").Value = "129.98"
$ws.Range("E6").Value = "  -1.53%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "0.577"
$ws.Range("E8").Value = "  -1.98%  "

$ws.Range("D9").Value = "2.316.97"
$ws.Range("E9").Value = "  -1.40%  "

$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("D11").Value = "5.58"
$ws.Range("E11").Value = "  +1.28%  "

$ws.Range("E12").Value = "  -0.57%  "

$ws.Range("D13").Value = "0.335"
$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").Value = "23.52"
$ws.Range("E14").Value = "  -1.54%  "

$ws.Range("D15").Value = "60.177.94"
$ws.Range("E15").Value = "  +0.09%  "

$ws.Range("D16").Value = "2.727.49"
$ws.Range("E16").Value = "  -1.75%  "

$ws.Range("E17").Value = "  +0.53%  "

$ws.Range("D18").Value = "2.308.30"
$ws.Range("E18").Value = "  -2.01%  "

$ws.Range("D19").Value = "10.56"
$ws.Range("E19").Value = "  -1.32%  "

$ws.Range("D20").Value = "4.08"
$ws.Range("E20").Value = "  -1.94%  "

$ws.Range("D21").Value = "313.33"
$ws.Range("E21").Value = "  -0.35%  "

$ws.Range("D22").Value = "6.56"
$ws.Range("E22").Value = "  -3.68%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").Value = "64.30"
$ws.Range("E24").Value = "  +1.57%  "

$ws.Range("D25").Value = "0.171"
$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").Value = "7.82"
$ws.Range("E27").Value = "  -0.81%  "

$ws.Range("E28").Value = "  +0.53%  "

$ws.Range("D29").Value = "1.26"
$ws.Range("E29").Value = "  +9.18%  "

$ws.Range("D30").Value = "171.31"
$ws.Range("E30").Value = "  -0.29%  "

$ws.Range("E31").Value = "  -1.11%  "

$ws.Range("D32").Value = "0.0₃0727"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").Value = "6.01"
$ws.Range("E33").Value = "  +1.16%  "

$ws.Range("B34").Value = "PolygonEcosystemToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D34").Value = "0.381"
$ws.Range("E34").Value = "  -0.34%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.35"
$ws.Range("E35").Value = "  -4.11%  "

$ws.Range("D36").Value = "17.96"
$ws.Range("E36").Value = "  -0.59%  "

$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.17%  "

$ws.Range("E39").Value = "  -2.14%  "

$ws.Range("D40").Value = "316.78"
$ws.Range("E40").Value = "  -0.95%  "

$ws.Range("D41").Value = "37.99"
$ws.Range("E41").Value = "  -0.35%  "

$ws.Range("D42").Value = "1.52"
$ws.Range("E42").Value = "  -0.96%  "

$ws.Range("D43").Value = "137.70"
$ws.Range("E43").Value = "  -3.13%  "

$ws.Range("D44").Value = "3.50"
$ws.Range("E44").Value = "  +0.97%  "

$ws.Range("D45").Value = "0.0943"
$ws.Range("E45").Value = "  -0.75%  "

$ws.Range("D46").Value = "19.07"
$ws.Range("E46").Value = "  -1.78%  "

$ws.Range("D47").Value = "0.566"
$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("E48").Value = "  -1.01%  "

$ws.Range("E49").Value = "  +0.26%  "

$ws.Range("E50").Value = "  +1.31%  "

$ws.Range("D51").Value = "10.92"
$ws.Range("E51").Value = "  -0.96%  "

$ws.Columns("D").NumberFormat = $dFormat
